$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBooks = @(
    '----------'
    'Extreme Sports'
    'Little Bear'
    'Earth Science and Weather'
    'Crocodile Encounters'
    'Grimmtastic Girls Sleeping Beauty Dreams Big'
    'A Bear Called Paddington'
    'Grimmtastic Girls Sleeping Beauty Dreams Big'
    'Giant Dinosaurs'
    'Calendar Mysteries February Friend'
    'The Yellow House Mystery'
    'A Den A Tree A Nest Is Best'
    'The Underground Railroad'
    'Nancy Drew The Secret Lost At Sea'
    'A Gift Of Heritage Historic Blacks in the Arts'
    'A Gift Of Heritage Historic Black Firsts'
    'A Wind in the Willows'
    'Philip Hall Likes Me I Reckon Maybe'
    'Aliens Dont Carve Jack O Lanterns'
    'The Giraffe and the Pelly and Me'
    'The Election Day Disaster'
    'Pete the Cat And His Four Groovy Buttons'
    'The Land of Many Colors'
    'Exploring Planets'
    'Focus on Water'
    'Chito Rodriguez-Gomez'
    'ZooBooks Wolves'
    'Thirteen Moons On Turtles Back'
    'Cordelia'
    'Shapes in the Sky'
    'A Rainbow All Around Me'
    'Exploring Weather'
    'If You Were A Minute'
    'Flakes and Fluries'
)

$startRow = 107
for ($i = 0; $i -lt $newBooks.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newBooks[$i]
}

$ws.Range("A$startRow").Select()
$excel.ActiveWindow.ScrollRow = $startRow
